$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Update the two teacher names (shared strings) used throughout column F.
#    Every cell in column F (rows 2-48) holds either "Литвинов Юрий
#    Викторович" or "Кириленко Яков Александрович"; append the requested
#    academic title suffix to each occurrence, wherever it is used.
# ---------------------------------------------------------------------------
$oldName1 = "Литвинов Юрий Викторович"
$newName1 = "Литвинов Юрий Викторович, доцент"
$oldName2 = "Кириленко Яков Александрович"
$newName2 = "Кириленко Яков Александрович, старший преподаватель"

for ($r = 2; $r -le 48; $r++) {
    $cell = $ws.Cells.Item($r, 6)
    $val = $cell.Value2
    if ($val -eq $oldName1) {
        $cell.Value = $newName1
    } elseif ($val -eq $oldName2) {
        $cell.Value = $newName2
    }
}

# ---------------------------------------------------------------------------
# 2) Widen column F so the longer teacher names stay readable.
# ---------------------------------------------------------------------------
$ws.Columns.Item(6).ColumnWidth = 44.8

# ---------------------------------------------------------------------------
# 3) Append new (currently empty) rows 49-103 below the existing data.
#    Rows 49-65 only carry formatting in columns B, H, I (same style as the
#    rest of the table). Rows 66-103 additionally carry a formatted, empty
#    cell in column D that uses a distinct font (Roboto, black) on the same
#    fill as the header/banner style.
# ---------------------------------------------------------------------------

# Template cells already carrying the styles we want to replicate.
$ws.Range("B1").Copy()
$ws.Range("B49:B103").PasteSpecial(-4122)
$ws.Range("H49:H103").PasteSpecial(-4122)
$ws.Range("I49:I103").PasteSpecial(-4122)

$ws.Range("C2").Copy()
$ws.Range("D66:D103").PasteSpecial(-4122)
$ws.Range("D66:D103").Font.Name = "Roboto"
$ws.Range("D66:D103").Font.Color = 0

$excel.CutCopyMode = 0
